$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Replace the four numeric PostalCode values in AddCustomerTest with new
# text values (qw22 / 22q2 / 2wer / rew3) -- new shared strings.
$ws1.Range("C3").Value = "qw22"
$ws1.Range("C4").Value = "22q2"
$ws1.Range("C6").Value = "2wer"
$ws1.Range("C7").Value = "rew3"

# Apply a text number format across the used range so every cell in
# A1:D7 picks up the new cellXfs entry (numFmtId 49 / "@").
$ws1.Range("A1:D7").NumberFormat = "@"

# Move the active sheet/selection from OpenAccount (A4:A5) to
# AddCustomerTest, selecting G9 -- matches the updated sheetViews.
$ws1.Activate()
$ws1.Range("G9").Select()
